$wb = $excel.ActiveWorkbook

$wsKarsten = $wb.Worksheets.Item("Karsten")
$wsDouwe = $wb.Worksheets.Item("Douwe")

# --- Karsten sheet (sheet1) ---
$wsKarsten.Range("D3").Value = "Nog 3 locaties toegevoegd met wat sublocaties in die locaties"
$wsKarsten.Range("D4").Value = "Nog 3 locaties toegevoegd met wat sublocaties in die locaties"
$wsKarsten.Range("B5").Value = 0.46875
$wsKarsten.Range("D5").Value = "Laatste voorlogie locatie toegevoegd en begin gemaakt aan de code"

# Update the selection shown on this sheet (cosmetic, matches diff)
$wsKarsten.Range("D6").Select()

# --- Douwe sheet (sheet2) ---
$wsDouwe.Range("D3").Value = "Nog 3 locaties toegevoegd met wat sublocaties in die locaties"
$wsDouwe.Range("B4").Value = 0.46875
$wsDouwe.Range("D4").Value = "Laatste voorlogie locatie toegevoegd en begin gemaakt aan de code"
$wsDouwe.Range("A5").Value = 0.55208333333333337
$wsDouwe.Range("B5").Value = 0.66666666666666663

# Update the selection shown on this sheet (cosmetic, matches diff) - keep it the active tab
$wsDouwe.Range("B5").Select()
